# Auto-generated edit script applying the Halicarnassus_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across the 8
# per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 452.625
$ws.Range("I2").Value = 375.25
$ws.Range("J2").Value = 530
$ws.Range("K2").Value = 375.25
$ws.Range("L2").Value = 530
$ws.Range("M2").Value = -262.25
$ws.Range("N2").Value = -756
# Row 20
$ws.Range("H20").Value = 3450.3333
$ws.Range("I20").Value = 380.7143
$ws.Range("K20").Value = 380.7143
$ws.Range("M20").Value = -150.7143
# Row 35
$ws.Range("H35").Value = 3450.3333
$ws.Range("I35").Value = 380.7143
$ws.Range("K35").Value = 380.7143
$ws.Range("M35").Value = -1.71429999999998
# Row 39
$ws.Range("H39").Value = 240.5
$ws.Range("I39").Value = 46.9
$ws.Range("J39").Value = 724.5
$ws.Range("K39").Value = 140.7
$ws.Range("L39").Value = 2173.5
$ws.Range("M39").Value = 155.3
$ws.Range("N39").Value = -2765.5
# Row 43
$ws.Range("H43").Value = 3375
$ws.Range("I43").Value = 3750
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 3750
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -3681
$ws.Range("N43").Value = -3138
# Row 58
$ws.Range("H58").Value = 1386.7
$ws.Range("J58").Value = 2700
$ws.Range("L58").Value = 8100
$ws.Range("N58").Value = -8400
# Row 61
$ws.Range("H61").Value = 8528.777
$ws.Range("I61").Value = 10108.429
$ws.Range("K61").Value = 30325.287
$ws.Range("M61").Value = -30153.287
# Row 80
$ws.Range("H80").Value = 931.5
$ws.Range("J80").Value = 1050
$ws.Range("L80").Value = 3150
$ws.Range("N80").Value = -5146
# Row 83
$ws.Range("H83").Value = 931.5
$ws.Range("J83").Value = 1050
$ws.Range("L83").Value = 9450
$ws.Range("N83").Value = -19434

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 1128.875
$ws.Range("I26").Value = 1128.875
$ws.Range("K26").Value = 1128.875
$ws.Range("M26").Value = -798.875
# Row 27
$ws.Range("H27").Value = 6099.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 6099.5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 6099.5
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -6467.5
# Row 74
$ws.Range("H74").Value = 5475.6665
$ws.Range("I74").Value = 5231.1333
$ws.Range("K74").Value = 5231.1333
$ws.Range("M74").Value = -4357.1333
# Row 77
$ws.Range("H77").Value = 5475.6665
$ws.Range("I77").Value = 5231.1333
$ws.Range("K77").Value = 26155.6665
$ws.Range("M77").Value = -21787.6665

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 5223.75
$ws.Range("I107").Value = 1385
$ws.Range("K107").Value = 1385
$ws.Range("M107").Value = 535
# Row 134
$ws.Range("H134").Value = 2586.6
$ws.Range("I134").Value = 2586.6
$ws.Range("K134").Value = 7759.799999999999
$ws.Range("M134").Value = -5224.799999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 2659
$ws.Range("I132").Value = 2487.7778
$ws.Range("K132").Value = 7463.3334
$ws.Range("M132").Value = -4933.3334
# Row 141
$ws.Range("H141").Value = 89270.375
$ws.Range("J141").Value = 89270.375
$ws.Range("L141").Value = 89270.375
$ws.Range("N141").Value = -99630.375

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 120000
$ws.Range("J37").Value = 120000
$ws.Range("L37").Value = 360000
$ws.Range("N37").Value = -360224
# Row 68
$ws.Range("H68").Value = 1099.5
$ws.Range("I68").Value = 1000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2189
# Row 71
$ws.Range("H71").Value = 1099.5
$ws.Range("I71").Value = 1000
$ws.Range("K71").Value = 9000
$ws.Range("M71").Value = -4944

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 3000
$ws.Range("J10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3338
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 113
$ws.Range("H113").Value = 5964.684
$ws.Range("I113").Value = 3332
$ws.Range("J113").Value = 8889.888999999999
$ws.Range("K113").Value = 3332
$ws.Range("L113").Value = 8889.888999999999
$ws.Range("M113").Value = -1162
$ws.Range("N113").Value = -13229.889
# Row 126
$ws.Range("H126").Value = 3374.75
$ws.Range("I126").Value = 1749.5
$ws.Range("K126").Value = 5248.5
$ws.Range("M126").Value = -2778.5
# Row 132
$ws.Range("H132").Value = 116767.89
$ws.Range("I132").Value = 205002.6
$ws.Range("K132").Value = 615007.8
$ws.Range("M132").Value = -612477.8

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3728.85
$ws.Range("I46").Value = 3631.75
$ws.Range("J46").Value = 3753.125
$ws.Range("K46").Value = 3631.75
$ws.Range("L46").Value = 3753.125
$ws.Range("M46").Value = -3443.75
$ws.Range("N46").Value = -4129.125
# Row 122
$ws.Range("H122").Value = 3473.75
$ws.Range("I122").Value = 2447.5
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 7342.5
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -4892.5
$ws.Range("N122").Value = -18400

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 81
$ws.Range("H81").Value = 790
$ws.Range("I81").Value = 737.5
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 1475
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -414
$ws.Range("N81").Value = -4122
# Row 84
$ws.Range("H84").Value = 790
$ws.Range("I84").Value = 737.5
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 7375
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -2071
$ws.Range("N84").Value = -20608
# Row 122
$ws.Range("H122").Value = 1647
$ws.Range("I122").Value = 1647
$ws.Range("K122").Value = 4941
$ws.Range("M122").Value = -2491
# Row 136
$ws.Range("H136").Value = 4435.6665
$ws.Range("I136").Value = 3990.75
$ws.Range("K136").Value = 11972.25
$ws.Range("M136").Value = -9422.25
